# Update TPM-derived NATMI edge statistics for Efna5-Epha4 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.1728506666666667"
$ws.Range("H2").Value = [double]"0.518552"
$ws.Range("I2").Value = [double]"0.0840503369699626"
$ws.Range("J2").Value = [double]"0.0840503369699626"
$ws.Range("M2").Value = [double]"8.813278666666667"
$ws.Range("N2").Value = [double]"26.439836"
$ws.Range("O2").Value = [double]"0.3770976991891536"
$ws.Range("P2").Value = [double]"0.3770976991891536"
$ws.Range("Q2").Value = [double]"1.523381093052445"
$ws.Range("R2").Value = [double]"13.710429837472"
$ws.Range("S2").Value = [double]"0.03169518868744595"
$ws.Range("T2").Value = [double]"0.03169518868744595"
$ws.Range("G3").Value = [double]"0.1728506666666667"
$ws.Range("H3").Value = [double]"0.518552"
$ws.Range("I3").Value = [double]"0.0840503369699626"
$ws.Range("J3").Value = [double]"0.0840503369699626"
$ws.Range("O3").Value = [double]"0.5522024902836482"
$ws.Range("P3").Value = [double]"0.5522024902836482"
$ws.Range("Q3").Value = [double]"2.230760980624889"
$ws.Range("R3").Value = [double]"20.076848825624"
$ws.Range("S3").Value = [double]"0.04641280538399313"
$ws.Range("T3").Value = [double]"0.04641280538399313"
$ws.Range("G4").Value = [double]"0.1728506666666667"
$ws.Range("H4").Value = [double]"0.518552"
$ws.Range("I4").Value = [double]"0.0840503369699626"
$ws.Range("J4").Value = [double]"0.0840503369699626"
$ws.Range("M4").Value = [double]"1.649921333333333"
$ws.Range("N4").Value = [double]"4.949764"
$ws.Range("O4").Value = [double]"0.07059592260441032"
$ws.Range("P4").Value = [double]"0.07059592260441033"
$ws.Range("Q4").Value = [double]"0.2851900024142223"
$ws.Range("R4").Value = [double]"2.566710021728"
$ws.Range("S4").Value = [double]"0.005933611083606087"
$ws.Range("T4").Value = [double]"0.005933611083606088"
$ws.Range("G5").Value = [double]"0.1728506666666667"
$ws.Range("H5").Value = [double]"0.518552"
$ws.Range("I5").Value = [double]"0.0840503369699626"
$ws.Range("J5").Value = [double]"0.0840503369699626"
$ws.Range("M5").Value = [double]"0.002428"
$ws.Range("N5").Value = [double]"0.007284"
$ws.Range("O5").Value = [double]"0.0001038879227879399"
$ws.Range("P5").Value = [double]"0.0001038879227879399"
$ws.Range("Q5").Value = [double]"0.0004196814186666667"
$ws.Range("R5").Value = [double]"0.003777132768"
$ws.Range("S5").Value = [double]"8.731814917435808e-06"
$ws.Range("T5").Value = [double]"8.731814917435808e-06"
$ws.Range("I6").Value = [double]"0.6650661694281633"
$ws.Range("J6").Value = [double]"0.6650661694281633"
$ws.Range("M6").Value = [double]"8.813278666666667"
$ws.Range("N6").Value = [double]"26.439836"
$ws.Range("O6").Value = [double]"0.3770976991891536"
$ws.Range("P6").Value = [double]"0.3770976991891536"
$ws.Range("Q6").Value = [double]"12.05407693365644"
$ws.Range("R6").Value = [double]"108.486692402908"
$ws.Range("S6").Value = [double]"0.2507949222999042"
$ws.Range("T6").Value = [double]"0.2507949222999042"
$ws.Range("I7").Value = [double]"0.6650661694281633"
$ws.Range("J7").Value = [double]"0.6650661694281633"
$ws.Range("O7").Value = [double]"0.5522024902836482"
$ws.Range("P7").Value = [double]"0.5522024902836482"
$ws.Range("S7").Value = [double]"0.3672511949616385"
$ws.Range("T7").Value = [double]"0.3672511949616385"
$ws.Range("I8").Value = [double]"0.6650661694281633"
$ws.Range("J8").Value = [double]"0.6650661694281633"
$ws.Range("M8").Value = [double]"1.649921333333333"
$ws.Range("N8").Value = [double]"4.949764"
$ws.Range("O8").Value = [double]"0.07059592260441032"
$ws.Range("P8").Value = [double]"0.07059592260441033"
$ws.Range("Q8").Value = [double]"2.256626556210222"
$ws.Range("R8").Value = [double]"20.309639005892"
$ws.Range("S8").Value = [double]"0.04695095982376226"
$ws.Range("T8").Value = [double]"0.04695095982376227"
$ws.Range("I9").Value = [double]"0.6650661694281633"
$ws.Range("J9").Value = [double]"0.6650661694281633"
$ws.Range("M9").Value = [double]"0.002428"
$ws.Range("N9").Value = [double]"0.007284"
$ws.Range("O9").Value = [double]"0.0001038879227879399"
$ws.Range("P9").Value = [double]"0.0001038879227879399"
$ws.Range("Q9").Value = [double]"0.003320818494666667"
$ws.Range("R9").Value = [double]"0.029887366452"
$ws.Range("S9").Value = [double]"6.909234285842402e-05"
$ws.Range("T9").Value = [double]"6.909234285842402e-05"
$ws.Range("G10").Value = [double]"0.5159453333333334"
$ws.Range("H10").Value = [double]"1.547836"
$ws.Range("I10").Value = [double]"0.2508834936018741"
$ws.Range("J10").Value = [double]"0.2508834936018741"
$ws.Range("M10").Value = [double]"8.813278666666667"
$ws.Range("N10").Value = [double]"26.439836"
$ws.Range("O10").Value = [double]"0.3770976991891536"
$ws.Range("P10").Value = [double]"0.3770976991891536"
$ws.Range("Q10").Value = [double]"4.54716999943289"
$ws.Range("R10").Value = [double]"40.924529994896"
$ws.Range("S10").Value = [double]"0.09460758820180348"
$ws.Range("T10").Value = [double]"0.09460758820180348"
$ws.Range("G11").Value = [double]"0.5159453333333334"
$ws.Range("H11").Value = [double]"1.547836"
$ws.Range("I11").Value = [double]"0.2508834936018741"
$ws.Range("J11").Value = [double]"0.2508834936018741"
$ws.Range("O11").Value = [double]"0.5522024902836482"
$ws.Range("P11").Value = [double]"0.5522024902836482"
$ws.Range("Q11").Value = [double]"6.658642051725779"
$ws.Range("R11").Value = [double]"59.927778465532"
$ws.Range("S11").Value = [double]"0.1385384899380166"
$ws.Range("T11").Value = [double]"0.1385384899380166"
$ws.Range("G12").Value = [double]"0.5159453333333334"
$ws.Range("H12").Value = [double]"1.547836"
$ws.Range("I12").Value = [double]"0.2508834936018741"
$ws.Range("J12").Value = [double]"0.2508834936018741"
$ws.Range("M12").Value = [double]"1.649921333333333"
$ws.Range("N12").Value = [double]"4.949764"
$ws.Range("O12").Value = [double]"0.07059592260441032"
$ws.Range("P12").Value = [double]"0.07059592260441033"
$ws.Range("Q12").Value = [double]"0.8512692123004445"
$ws.Range("R12").Value = [double]"7.661422910704"
$ws.Range("S12").Value = [double]"0.01771135169704198"
$ws.Range("T12").Value = [double]"0.01771135169704198"
$ws.Range("G13").Value = [double]"0.5159453333333334"
$ws.Range("H13").Value = [double]"1.547836"
$ws.Range("I13").Value = [double]"0.2508834936018741"
$ws.Range("J13").Value = [double]"0.2508834936018741"
$ws.Range("M13").Value = [double]"0.002428"
$ws.Range("N13").Value = [double]"0.007284"
$ws.Range("O13").Value = [double]"0.0001038879227879399"
$ws.Range("P13").Value = [double]"0.0001038879227879399"
$ws.Range("Q13").Value = [double]"0.001252715269333333"
$ws.Range("R13").Value = [double]"0.011274437424"
$ws.Range("S13").Value = [double]"2.606376501208012e-05"
$ws.Range("T13").Value = [double]"2.606376501208012e-05"
